$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This script applies the IFRS restatement edit: columns D-AJ for rows 2-6 are
# rewritten with corrected (much smaller / differently-signed) figures, columns
# J ("당기순이익(비지배)") and O ("자본총계(비지배)") are cleared for rows 2-6,
# and rows 7-9 have all of their financial data (columns D through AJ) cleared,
# leaving only the row index (A) and label columns (B, C) intact.

# Row 2
$ws.Range("D2").Value = 855
$ws.Range("E2").Value = -86
$ws.Range("F2").Value = -86
$ws.Range("G2").Value = -102
$ws.Range("H2").Value = -112
$ws.Range("I2").Value = -112
$ws.Range("K2").Value = 1083
$ws.Range("L2").Value = 529
$ws.Range("M2").Value = 554
$ws.Range("N2").Value = 554
$ws.Range("P2").Value = 55
$ws.Range("Q2").Value = -119
$ws.Range("R2").Value = -175
$ws.Range("S2").Value = 73
$ws.Range("T2").Value = 53
$ws.Range("U2").Value = -172
$ws.Range("V2").Value = 362
$ws.Range("W2").Value = -10.08
$ws.Range("X2").Value = -13.09
$ws.Range("Y2").Value = -17.84
$ws.Range("Z2").Value = -9.99
$ws.Range("AA2").Value = 95.42
$ws.Range("AB2").Value = 959.14
$ws.Range("AC2").Value = -1957
$ws.Range("AD2").Value = -2.73
$ws.Range("AE2").Value = 10412
$ws.Range("AF2").Value = 0.51
$ws.Range("AG2").Value = 192
$ws.Range("AH2").Value = 3.6
$ws.Range("AI2").Value = -9.119999999999999
$ws.Range("AJ2").Value = 5718505
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 881
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 16
$ws.Range("G3").Value = 11
$ws.Range("H3").Value = 6
$ws.Range("I3").Value = 6
$ws.Range("K3").Value = 993
$ws.Range("L3").Value = 449
$ws.Range("M3").Value = 544
$ws.Range("N3").Value = 544
$ws.Range("P3").Value = 55
$ws.Range("Q3").Value = 29
$ws.Range("R3").Value = 94
$ws.Range("S3").Value = -102
$ws.Range("T3").Value = 24
$ws.Range("U3").Value = 6
$ws.Range("V3").Value = 282
$ws.Range("W3").Value = 1.82
$ws.Range("X3").Value = 0.67
$ws.Range("Y3").Value = 1.07
$ws.Range("Z3").Value = 0.57
$ws.Range("AA3").Value = 82.52
$ws.Range("AB3").Value = 935.92
$ws.Range("AC3").Value = 103
$ws.Range("AD3").Value = 67.58
$ws.Range("AE3").Value = 10224
$ws.Range("AF3").Value = 0.68
$ws.Range("AG3").Value = 192
$ws.Range("AH3").Value = 2.77
$ws.Range("AI3").Value = 173.48
$ws.Range("AJ3").Value = 5718505
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 968
$ws.Range("E4").Value = 39
$ws.Range("F4").Value = 39
$ws.Range("G4").Value = 27
$ws.Range("H4").Value = 9
$ws.Range("I4").Value = 9
$ws.Range("K4").Value = 1010
$ws.Range("L4").Value = 481
$ws.Range("M4").Value = 529
$ws.Range("N4").Value = 529
$ws.Range("P4").Value = 55
$ws.Range("Q4").Value = 56
$ws.Range("R4").Value = -29
$ws.Range("S4").Value = -1
$ws.Range("T4").Value = 41
$ws.Range("U4").Value = 16
$ws.Range("V4").Value = 302
$ws.Range("W4").Value = 3.98
$ws.Range("X4").Value = 0.92
$ws.Range("Y4").Value = 1.65
$ws.Range("Z4").Value = 0.89
$ws.Range("AA4").Value = 91.02
$ws.Range("AB4").Value = 907.9299999999999
$ws.Range("AC4").Value = 155
$ws.Range("AD4").Value = 51.88
$ws.Range("AE4").Value = 9934
$ws.Range("AF4").Value = 0.8100000000000001
$ws.Range("AG4").Value = 144
$ws.Range("AH4").Value = 1.79
$ws.Range("AI4").Value = 86.34
$ws.Range("AJ4").Value = 5718505
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 920
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = 13
$ws.Range("G5").Value = -7
$ws.Range("H5").Value = -13
$ws.Range("I5").Value = -13
$ws.Range("K5").Value = 1355
$ws.Range("L5").Value = 843
$ws.Range("M5").Value = 512
$ws.Range("N5").Value = 512
$ws.Range("P5").Value = 55
$ws.Range("Q5").Value = 38
$ws.Range("R5").Value = -375
$ws.Range("S5").Value = 285
$ws.Range("T5").Value = 362
$ws.Range("U5").Value = -324
$ws.Range("V5").Value = 612
$ws.Range("W5").Value = 1.42
$ws.Range("X5").Value = -1.38
$ws.Range("Y5").Value = -2.44
$ws.Range("Z5").Value = -1.07
$ws.Range("AA5").Value = 164.77
$ws.Range("AB5").Value = 876.96
$ws.Range("AC5").Value = -221
$ws.Range("AD5").Value = -36.83
$ws.Range("AE5").Value = 9614
$ws.Range("AF5").Value = 0.85
$ws.Range("AG5").Value = 144
$ws.Range("AH5").Value = 1.77
$ws.Range("AI5").Value = -60.44
$ws.Range("AJ5").Value = 5718505
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 946
$ws.Range("E6").Value = -57
$ws.Range("F6").Value = -57
$ws.Range("G6").Value = -81
$ws.Range("H6").Value = -86
$ws.Range("I6").Value = -86
$ws.Range("K6").Value = 1479
$ws.Range("L6").Value = 919
$ws.Range("M6").Value = 561
$ws.Range("N6").Value = 561
$ws.Range("P6").Value = 65
$ws.Range("Q6").Value = -26
$ws.Range("R6").Value = -148
$ws.Range("S6").Value = 195
$ws.Range("T6").Value = 37
$ws.Range("U6").Value = -63
$ws.Range("V6").Value = 667
$ws.Range("W6").Value = -6.06
$ws.Range("X6").Value = -9.050000000000001
$ws.Range("Y6").Value = -15.97
$ws.Range("Z6").Value = -6.04
$ws.Range("AA6").Value = 163.8
$ws.Range("AB6").Value = 806.86
$ws.Range("AC6").Value = -1413
$ws.Range("AD6").Value = -14.58
$ws.Range("AE6").Value = 9189
$ws.Range("AF6").Value = 2.24
$ws.Range("AG6").Value = 100
$ws.Range("AH6").Value = 0.49
$ws.Range("AI6").Value = -7.13
$ws.Range("AJ6").Value = 6500000

# Rows 7-9: clear all data columns (D:AJ), keep only A (index) and C (label)
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
